$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("problem_input")
$wsRobot = $wb.Worksheets.Item("robot")

# --- "problem_input" sheet: column A header changes from the old
# "average service time per parcel (seconds)" label to "postcode". ---
$wsInput.Range("A1").Value = "postcode"

# --- "robot" sheet: C1/C2 now describe the service rate instead of the
# (now unused) "postcode" label that used to live there. Do this before
# adding the new header on "problem_input" below so the shared-string
# table is built up in the same order as the target workbook (service
# rate before customer arrival rate). ---
$wsRobot.Range("C1").Value = "service rate (per minute)"
$wsRobot.Range("C2").Value = 5

# --- "problem_input" sheet: new column E with header + per-row data. ---
$eHeader = $wsInput.Range("E1")
$eHeader.Value = "customer arrival rate (per minute)"
$eHeader.Font.Bold = $true
$eHeader.HorizontalAlignment = -4131   # xlLeft
$eHeader.VerticalAlignment = -4160     # xlTop
$eHeader.Borders.LineStyle = 1         # xlContinuous
$eHeader.Borders.Weight = 2            # xlThin

$wsInput.Range("E2").Value = 0
$wsInput.Range("E3").Value = 5
$wsInput.Range("E4").Value = 3
$wsInput.Range("E5").Value = 5
$wsInput.Range("E6").Value = 3
$wsInput.Range("E7").Value = 5

# --- Update the stored selections on each sheet to match the edited
# workbook, while leaving "problem_input" as the active tab. ---
$wsRobot.Range("C3").Select()
$wsInput.Activate()
$wsInput.Range("E7").Select()
